# Insert a new data row at row 43 (weekly update for
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Repollo").
# This pushes the existing rows 43-97 down to 44-98 and extends
# the sheet's used range to A1:R98.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(43).Insert()

$ws.Range("A43").Value() = 7
$ws.Range("B43").Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C43").Value() = "Ñuble"
$ws.Range("D43").Value() = 44413
$ws.Range("E43").Value() = 16
$ws.Range("F43").Value() = 100112006
$ws.Range("G43").Value() = "Repollo"
$ws.Range("H43").Value() = "Crespo record"
$ws.Range("I43").Value() = "Primera"
$ws.Range("J43").Value() = 120
$ws.Range("K43").Value() = 600
$ws.Range("L43").Value() = 650
$ws.Range("M43").Value() = 625
$ws.Range("N43").Value() = "$/unidad"
$ws.Range("O43").Value() = "Provincia de Diguillín"
$ws.Range("P43").Value() = 625
$ws.Range("Q43").Value() = 1
$ws.Range("R43").Value() = "Hortaliza"
